$wb = $excel.ActiveWorkbook

# Overview sheet: row for 9a935e24-... (row 3) moves from
# "Ready for handoff" to "Handed back: in sync with en-US" for both locales.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: Status -> Handed back, and record the handback datetime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-12 14:37:43"

# de-de detail sheet: Status -> Handed back, and record the handback datetime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-12 14:37:49"
